# Update the betting odds for row 2 (the single match row) to reflect
# the latest FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 2.22
    "H2"  = 3.25
    "I2"  = 3.15
    "J2"  = 2.75
    "K2"  = 2.12
    "S2"  = 1.37
    "T2"  = 2.85
    "U2"  = 1.6
    "V2"  = 2.22
    "W2"  = 9
    "X2"  = 12
    "Z2"  = 23
    "AB2" = 22
    "AE2" = 11.75
    "AH2" = 10.75
    "AI2" = 17.5
    "AK2" = 40
    "AM2" = 29
    "AN2" = 4.3
    "AO2" = 11.5
    "AP2" = 17
    "AR2" = 65
    "AS2" = 175
    "AT2" = 2.85
    "AU2" = 6.5
    "AW2" = 5.2
    "AY2" = 22
    "BA2" = 110
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
